$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.23229999999999
$ws.Range("C14").Value = -11.9797
$ws.Range("C21").Value = -13.04900000000001
$ws.Range("C23").Value = -11.8755
$ws.Range("C25").Value = -11.26009999999999
